# Add season-record columns (Wins / Losses / Ties) to the stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: copy the formatting of the last existing header cell (AC1, style s="1")
# onto the three new header cells, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "Wins"

$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "Losses"

$ws.Range("AC1").Copy($ws.Range("AF1"))
$ws.Range("AF1").Value = "Ties"

# Data rows 2-46: every player on the roster shares the team's 2009 season
# record (88 wins, 74 losses, 0 ties).
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}
